# Add the season-record columns ("Wins", "Losses", "Ties") to the
# player-stats table on Sheet1, matching the format of the existing
# header columns and filling every player's row with the team's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered style used by the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill the season record for every player row (2-46).
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 75   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 87   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
